$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("G12").Value = 1240524717.0500007
$ws.Range("I12").Value = 4188377156

# Row 13
$ws.Range("G13").Value = 319819483.18000001
$ws.Range("I13").Value = 1012006300

# Row 14
$ws.Range("G14").Value = 34063116.800000042
$ws.Range("I14").Value = -44319159.289999999

# Row 16
$ws.Range("I16").Value = -162861893.59999999

# Row 18 - G18 becomes a formula (was a hard-coded value before)
$ws.Range("G18").Formula = "=SUM(G12:G17)"

# Row 19 - fill in previously blank G19, flip sign of I19
$ws.Range("G19").Value = -379300000.00000012
$ws.Range("I19").Value = 1160500000

# Row 21 - fill in previously blank G21 with a formula
$ws.Range("G21").Formula = "=SUM(G18:G20)"

# Row 22 - fill in previously blank G22
$ws.Range("G22").Value = -20015625

# Row 26 - G26 was a shared-string placeholder; now a numeric value
$ws.Range("G26").Value = 1029174575.116062
$ws.Range("I26").Value = 1010658959
